$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2755.75
$ws.Range("I40").Value = 3666
$ws.Range("K40").Value = 3666
$ws.Range("M40").Value = -3491
$ws.Range("H45").Value = 808
$ws.Range("I45").Value = 808
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2424
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2232
$ws.Range("N45").ClearContents()
$ws.Range("H98").Value = 3304
$ws.Range("I98").Value = 2111.4285
$ws.Range("K98").Value = 2111.4285
$ws.Range("M98").Value = -613.4285
$ws.Range("H122").Value = 3304
$ws.Range("I122").Value = 2111.4285
$ws.Range("K122").Value = 6334.2855
$ws.Range("M122").Value = -3884.2855
$ws.Range("H129").Value = 100002140
$ws.Range("I129").Value = 111111520
$ws.Range("K129").Value = 333334560
$ws.Range("M129").Value = -333329560
$ws.Range("H132").Value = 1736.3733
$ws.Range("I132").Value = 1714.7747
$ws.Range("K132").Value = 5144.3241
$ws.Range("M132").Value = -2614.3241
$ws.Range("H137").Value = 2542.1943
$ws.Range("I137").Value = 2431.8262
$ws.Range("J137").Value = 2737.4614
$ws.Range("K137").Value = 7295.4786
$ws.Range("L137").Value = 8212.3842
$ws.Range("M137").Value = -4745.4786
$ws.Range("N137").Value = -13312.3842
$ws.Range("H138").Value = 4877.8486
$ws.Range("I138").Value = 2736
$ws.Range("J138").Value = 5145.5796
$ws.Range("K138").Value = 8208
$ws.Range("L138").Value = 15436.7388
$ws.Range("M138").Value = -3068
$ws.Range("N138").Value = -25716.7388

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12378.638
$ws.Range("I32").Value = 9591.370000000001
$ws.Range("J32").Value = 50006.75
$ws.Range("K32").Value = 9591.370000000001
$ws.Range("L32").Value = 50006.75
$ws.Range("M32").Value = -9304.370000000001
$ws.Range("N32").Value = -50580.75
$ws.Range("H122").Value = 1588.1111
$ws.Range("I122").Value = 1570.2174
$ws.Range("K122").Value = 4710.6522
$ws.Range("M122").Value = -2260.6522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 55816.332
$ws.Range("J59").Value = 55816.332
$ws.Range("L59").Value = 55816.332
$ws.Range("N59").Value = -57510.332
$ws.Range("H86").Value = 2231.5
$ws.Range("I86").Value = 2198.9
$ws.Range("K86").Value = 2198.9
$ws.Range("M86").Value = -1075.9
$ws.Range("H89").Value = 2231.5
$ws.Range("I89").Value = 2198.9
$ws.Range("K89").Value = 10994.5
$ws.Range("M89").Value = -5378.5
$ws.Range("H94").Value = 592.5333000000001
$ws.Range("I94").Value = 659
$ws.Range("K94").Value = 659
$ws.Range("M94").Value = -208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5481.364
$ws.Range("I16").Value = 4702.5
$ws.Range("K16").Value = 4702.5
$ws.Range("M16").Value = -4415.5
$ws.Range("H31").Value = 3497.476
$ws.Range("I31").Value = 3161.8235
$ws.Range("J31").Value = 4924
$ws.Range("K31").Value = 3161.8235
$ws.Range("L31").Value = 4924
$ws.Range("M31").Value = -2866.8235
$ws.Range("N31").Value = -5514
$ws.Range("H34").Value = 3497.476
$ws.Range("I34").Value = 3161.8235
$ws.Range("J34").Value = 4924
$ws.Range("K34").Value = 3161.8235
$ws.Range("L34").Value = 4924
$ws.Range("M34").Value = -2959.8235
$ws.Range("N34").Value = -5328
$ws.Range("H58").Value = 3362.9285
$ws.Range("I58").Value = 3379.8
$ws.Range("J58").Value = 3353.5557
$ws.Range("K58").Value = 3379.8
$ws.Range("L58").Value = 3353.5557
$ws.Range("M58").Value = -3176.8
$ws.Range("N58").Value = -3759.5557
$ws.Range("H105").Value = 3065.5557
$ws.Range("I105").Value = 2942.7144
$ws.Range("K105").Value = 2942.7144
$ws.Range("M105").Value = -1195.7144
$ws.Range("H113").Value = 5481.364
$ws.Range("I113").Value = 4702.5
$ws.Range("K113").Value = 4702.5
$ws.Range("M113").Value = -2532.5
$ws.Range("H134").Value = 24213.268
$ws.Range("I134").Value = 21466.666
$ws.Range("K134").Value = 64399.99800000001
$ws.Range("M134").Value = -61864.99800000001
$ws.Range("H136").Value = 3362.9285
$ws.Range("I136").Value = 3379.8
$ws.Range("J136").Value = 3353.5557
$ws.Range("K136").Value = 10139.4
$ws.Range("L136").Value = 10060.6671
$ws.Range("M136").Value = -7589.400000000001
$ws.Range("N136").Value = -15160.6671
$ws.Range("H141").Value = 265555.8
$ws.Range("J141").Value = 275362.7
$ws.Range("L141").Value = 275362.7
$ws.Range("N141").Value = -285722.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5132.5
$ws.Range("I3").Value = 1865
$ws.Range("J3").Value = 8400
$ws.Range("K3").Value = 5595
$ws.Range("L3").Value = 25200
$ws.Range("M3").Value = -5483
$ws.Range("N3").Value = -25424
$ws.Range("H58").Value = 1832.3334
$ws.Range("I58").Value = 1832.3334
$ws.Range("K58").Value = 5497.0002
$ws.Range("M58").Value = -5369.0002
$ws.Range("H107").Value = 309
$ws.Range("I107").Value = 368.6
$ws.Range("J107").Value = 209.66667
$ws.Range("K107").Value = 1105.8
$ws.Range("L107").Value = 629.00001
$ws.Range("M107").Value = 814.1999999999998
$ws.Range("N107").Value = -4469.00001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2562.75
$ws.Range("I97").Value = 3190.75
$ws.Range("K97").Value = 3190.75
$ws.Range("M97").Value = -2694.75
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 6997.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 6997.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 6997.5
$ws.Range("N3").Value = -7221.5
$ws.Range("M3").ClearContents()
$ws.Range("H15").Value = 6997.5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 6997.5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 6997.5
$ws.Range("N15").Value = -7337.5
$ws.Range("M15").ClearContents()
$ws.Range("H19").Value = 381
$ws.Range("I19").Value = 71.5
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 71.5
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = 98.5
$ws.Range("N19").Value = -1340
$ws.Range("H22").Value = 3907.8
$ws.Range("I22").Value = 2850.8572
$ws.Range("J22").Value = 4476.923
$ws.Range("K22").Value = 2850.8572
$ws.Range("L22").Value = 4476.923
$ws.Range("M22").Value = -2555.8572
$ws.Range("N22").Value = -5066.923
$ws.Range("H27").Value = 3907.8
$ws.Range("I27").Value = 2850.8572
$ws.Range("J27").Value = 4476.923
$ws.Range("K27").Value = 2850.8572
$ws.Range("L27").Value = 4476.923
$ws.Range("M27").Value = -2743.8572
$ws.Range("N27").Value = -4690.923
$ws.Range("H122").Value = 352672.7
$ws.Range("I122").Value = 507636.44
$ws.Range("K122").Value = 1522909.32
$ws.Range("M122").Value = -1520459.32
$ws.Range("H132").Value = 5219.68
$ws.Range("I132").Value = 4521.6
$ws.Range("J132").Value = 6266.8
$ws.Range("K132").Value = 13564.8
$ws.Range("L132").Value = 18800.4
$ws.Range("M132").Value = -11034.8
$ws.Range("N132").Value = -23860.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H62").Value = 17306.46
$ws.Range("I62").Value = 13284.143
$ws.Range("J62").Value = 21999.166
$ws.Range("K62").Value = 13284.143
$ws.Range("L62").Value = 21999.166
$ws.Range("M62").Value = -12660.143
$ws.Range("N62").Value = -23247.166
$ws.Range("H65").Value = 17306.46
$ws.Range("I65").Value = 13284.143
$ws.Range("J65").Value = 21999.166
$ws.Range("K65").Value = 66420.715
$ws.Range("L65").Value = 109995.83
$ws.Range("M65").Value = -63300.715
$ws.Range("N65").Value = -116235.83
$ws.Range("H100").Value = 631.1429000000001
$ws.Range("I100").Value = 279.5
$ws.Range("K100").Value = 559
$ws.Range("M100").Value = -18
$ws.Range("H107").Value = 709.8387
$ws.Range("I107").Value = 583.6667
$ws.Range("K107").Value = 1751.0001
$ws.Range("M107").Value = 168.9999
$ws.Range("H113").Value = 771600.3
$ws.Range("I113").Value = 1980.2
$ws.Range("K113").Value = 5940.6
$ws.Range("M113").Value = -3770.6
$ws.Range("H126").Value = 10626.692
$ws.Range("I126").Value = 11522
$ws.Range("J126").Value = 5702.5
$ws.Range("K126").Value = 34566
$ws.Range("L126").Value = 17107.5
$ws.Range("M126").Value = -32096
$ws.Range("N126").Value = -22047.5
$ws.Range("H132").Value = 4405.737
$ws.Range("J132").Value = 4067.7856
$ws.Range("L132").Value = 12203.3568
$ws.Range("N132").Value = -17263.3568
$ws.Range("H136").Value = 2708.04
$ws.Range("I136").Value = 2577.318
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 7731.954000000001
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -5181.954000000001
$ws.Range("N136").Value = -16100.0001
